# BKD_QTR_FIN.xlsx quarterly refresh
# Inserts two new fiscal-quarter columns (D:E) ahead of the existing data
# (old D:K shifts to F:M), fills the two new quarters with the newly
# reported figures, and applies a couple of restated historical figures
# that came back from the data provider together with this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two blank columns at D:E -------------------------
# Everything that used to live in D:K slides right to F:M.
$ws.Columns("D:E").Insert()

# Column F used to be column D, so it still carries the right number
# formats/styles (date format for the "Period Ending" rows, plain number
# format elsewhere). Clone that formatting into the two new columns so the
# new cells render the same way the rest of the table does.
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)

# --- 2. New-quarter data for columns D (latest) and E (prior quarter) ------
# Keyed by row number -> @(D-value, E-value)
$data = @{
    7   = @(43465, 43373)
    8   = @(1068900, 1120100)
    9   = @(831300, 868400)
    10  = @(237600, 251700)
    12  = @("NA", "NA")
    13  = @(0, 0)
    14  = @(62700, 7900)
    15  = @(106100, 111000)
    17  = @(1127000, 1116500)
    18  = @(-58100, 3600)
    20  = @(222800, 10100)
    21  = @(271500, 125500)
    22  = @(65000, 68600)
    23  = @(99800, -54900)
    24  = @(-37800, -17800)
    25  = @(0, 0)
    26  = @(137600, -37100)
    27  = @(137600, -37100)
    28  = @(0, 0)
    29  = @(-6000, "NA")
    30  = @(0, 0)
    31  = @(0, 0)
    32  = @(-222800, -10100)
    33  = @(131500, -37100)
    34  = @(0, 0)
    35  = @(131500, -37100)
    38  = @(43465, 43373)
    41  = @(398300, 133700)
    42  = @(14900, 0)
    43  = @(133900, 130100)
    44  = @(0, 0)
    45  = @(227000, 373000)
    46  = @(774000, 636700)
    47  = @(27500, 30000)
    48  = @(5275400, 5407100)
    49  = @(205600, 213800)
    50  = @(0, 0)
    51  = @(0, 0)
    52  = @(184700, 207000)
    53  = @(0, 0)
    54  = @(6467300, 6494600)
    57  = @(95000, 79600)
    58  = @(317600, 503700)
    59  = @(360700, 379500)
    60  = @(773300, 962800)
    61  = @(4197100, 4129300)
    62  = @(478400, 512700)
    63  = @(0, 0)
    64  = @(0, 0)
    65  = @(0, 0)
    66  = @(5448400, 5604300)
    68  = @(0, 0)
    69  = @(0, 0)
    70  = @(0, 0)
    71  = @(0, 0)
    72  = @(-3069300, -3200800)
    73  = @(0, 0)
    74  = @(0, 0)
    75  = @(0, 0)
    76  = @(1018900, 890400)
    77  = @(0, 0)
    80  = @(43465, 43373)
    81  = @(131500, -37100)
    83  = @(106700, 111800)
    84  = @(0, 0)
    85  = @(0, 0)
    86  = @(0, 0)
    87  = @(0, 0)
    88  = @(0, 0)
    89  = @(33500, 71900)
    91  = @(-56100, -48900)
    92  = @(0, 0)
    93  = @(0, 0)
    94  = @(301800, -24500)
    96  = @(0, 0)
    97  = @(0, 0)
    98  = @(0, 0)
    99  = @(0, 0)
    100 = @(-85100, -37900)
    101 = @(0, 0)
    102 = @(250100, 9400)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
}

# --- 3. A handful of historical cells were restated in this refresh --------
# (values shifted into columns F:J but don't match the old pre-shift figures)
$ws.Cells.Item(89, 8).Value = 95300

$ws.Cells.Item(91, 6).Value = -93100
$ws.Cells.Item(91, 7).Value = -27300
$ws.Cells.Item(91, 8).Value = -4800
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 48500

Write-Output "BKD quarterly refresh applied"
